$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "87.997.79"
$ws.Range("E2").Value2 = "  +7.76%  "
$ws.Range("D3").Value2 = "3.366.12"
$ws.Range("E3").Value2 = "  +6.31%  "
$ws.Range("E4").Value2 = "  +0.43%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "217.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  +3.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "641.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  +3.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.413"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value2 = "  +46.95%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.664"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value2 = "  +14.21%  "
$ws.Range("E9").Value2 = "  +0.19%  "
$ws.Range("D10").Value2 = "3.363.96"
$ws.Range("E10").Value2 = "  +6.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.614"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = "  +5.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.0000282"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = "  +12.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "36.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value2 = "  +15.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "0.166"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  +1.01%  "
$ws.Range("D15").Value2 = "4.008.62"
$ws.Range("E15").Value2 = "  +7.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "5.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = "  +3.61%  "
$ws.Range("D17").Value2 = "88.059.16"
$ws.Range("E17").Value2 = "  +8.38%  "
$ws.Range("D18").Value2 = "3.378.76"
$ws.Range("E18").Value2 = "  +6.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "14.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "  +6.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "9.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  +7.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "451.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "  +4.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "3.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "  -4.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "5.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "  +9.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "7.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "  +2.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "5.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = "  +4.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "12.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "  +14.95%  "
$ws.Range("D27").Value2 = "3.577.45"
$ws.Range("E27").Value2 = "  +8.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "80.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "  +5.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "0.0000141"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = "  +17.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "  -0.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "0.186"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "  +35.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "9.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = "  +3.91%  "
$ws.Range("E33").Value2 = "  +0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "568.25"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = "  -3.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "1.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = "  +1.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "2.08"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "  +4.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "7.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = "  +18.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "0.141"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = "  -8.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "23.61"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = "  +4.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "0.430"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = "  +5.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "21.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = "  +5.34%  "
$ws.Range("E42").Value2 = "  +0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "2.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = "  +1.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "3.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = "  +1.88%  "
$ws.Range("E45").Value2 = "  -0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "157.46"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "  -1.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "186.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  -0.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "1.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  +5.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "46.47"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = "  +3.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "4.46"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "  +6.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "0.661"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = "  +5.97%  "
